# Objectif 1 avec nvx deltas
$wb = $excel.ActiveWorkbook

# --- Sheet "Base Deltas" ---
$ws1 = $wb.Worksheets.Item("Base Deltas")
$ws1.Range("E4").Value = -0.22
$ws1.Range("F4").Value = -0.01
$ws1.Range("I4").Value = 0.05
$ws1.Range("G5").Value = 0.09
$ws1.Range("F6").Value = -0.62
$ws1.Range("F8").Value = -0.05
$ws1.Range("D9").Value = -0.04

# --- Sheet "Calculated Deltas" ---
$ws2 = $wb.Worksheets.Item("Calculated Deltas")
$ws2.Range("D4").Value = -0.23
$ws2.Range("E4").Value = 0.43
$ws2.Range("E5").Value = -0.7000000000000001
$ws2.Range("D6").Value = -1.2
$ws2.Range("E6").Value = -1.31
$ws2.Range("D8").Value = 0.28
$ws2.Range("E8").Value = 0.25
